$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update numeric columns B:K (col A text unchanged)
$ws.Cells.Item(2, 2).Value = 0.38221848321154955
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 0

# Row 3: update numeric columns B:K (col A text unchanged)
$ws.Cells.Item(3, 2).Value = 0.32602415205646584
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = -0.05136370790070153
$ws.Cells.Item(3, 5).Value = -0.002319931404911704
$ws.Cells.Item(3, 6).Value = 0.0026846829289998934
$ws.Cells.Item(3, 7).Value = 0.0017865144626340477
$ws.Cells.Item(3, 8).Value = 0.0011051303385501894
$ws.Cells.Item(3, 9).Value = -0.012547107046802803
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 0.0044600874671482105

# Row 4: update numeric columns B:K (col A text unchanged)
$ws.Cells.Item(4, 2).Value = 0.3134393300723387
$ws.Cells.Item(4, 3).Value = 0.009338056976012147
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 0.00036644361940457265
$ws.Cells.Item(4, 6).Value = 0.000425252531844555
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 0.001305368400552078
$ws.Cells.Item(4, 9).Value = -0.023742784732011146
$ws.Cells.Item(4, 10).Value = 0.00018045717558565327
$ws.Cells.Item(4, 11).Value = -0.00045761595551502277

# Row 5: update numeric columns B:K (col A text unchanged)
$ws.Cells.Item(5, 2).Value = 0.41778867121429536
$ws.Cells.Item(5, 3).Value = 0.14854865321703328
$ws.Cells.Item(5, 4).Value = -0.07104837408526352
$ws.Cells.Item(5, 5).Value = 0.006429121631376614
$ws.Cells.Item(5, 6).Value = 0.03838564389616107
$ws.Cells.Item(5, 7).Value = -0.0034918244950533917
$ws.Cells.Item(5, 8).Value = -0.00019718102002140152
$ws.Cells.Item(5, 9).Value = -0.012080462543896421
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = -0.002196235458379603

# Row 6: update numeric columns B:K (col A text unchanged)
$ws.Cells.Item(6, 2).Value = 0.49796928183883504
$ws.Cells.Item(6, 3).Value = 0.04702823849462218
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = -0.0028029053001516423
$ws.Cells.Item(6, 6).Value = 0.00443604947308156
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 0.007271269855284854
$ws.Cells.Item(6, 9).Value = 0.024027018456648928
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 0.00022093964505381924

# Row 7: update numeric columns B:K (col A text unchanged)
$ws.Cells.Item(7, 2).Value = 0.3827651650712497
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 0.08376162806595232
$ws.Cells.Item(7, 5).Value = -0.008367600269295231
$ws.Cells.Item(7, 6).Value = -0.23064411216691255
$ws.Cells.Item(7, 7).Value = 0.0011100858778281828
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0.009234229069589038
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 0.02970165265525293

# Row 8: update numeric columns B:K (col A text unchanged)
$ws.Cells.Item(8, 2).Value = -0.05152151303973196
$ws.Cells.Item(8, 3).Value = -0.4327120159910155
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 0.0007287597933152894
$ws.Cells.Item(8, 6).Value = 0.015484888633800674
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = -0.0002822923027926719
$ws.Cells.Item(8, 9).Value = -0.016701385632888577
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = -0.0008046326114007907

# Row 9: update numeric columns B:K (col A text unchanged)
$ws.Cells.Item(9, 2).Value = -0.076412523990948
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = -0.08725317411722248
$ws.Cells.Item(9, 5).Value = -0.0024345807838146226
$ws.Cells.Item(9, 6).Value = 0.07123614705340907
$ws.Cells.Item(9, 7).Value = -0.005718842696492503
$ws.Cells.Item(9, 8).Value = 0.00019360377175146762
$ws.Cells.Item(9, 9).Value = -0.0006496023095945419
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = -0.0002645618692524304

# Row 10: update numeric columns B:K (col A text unchanged)
$ws.Cells.Item(10, 2).Value = 0.18058404175502063
$ws.Cells.Item(10, 3).Value = 0.30381285795704704
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = 0.000009589721301508013
$ws.Cells.Item(10, 6).Value = -0.0013851645435718708
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = 0.0002653294047978285
$ws.Cells.Item(10, 9).Value = 0.01860021426120921
$ws.Cells.Item(10, 10).Value = -0.07346623000293778
$ws.Cells.Item(10, 11).Value = 0.009159968948122682

# Row 11: update numeric columns B:K (col A text unchanged)
$ws.Cells.Item(11, 2).Value = 0.5317091790399617
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = 0.25826017480032265
$ws.Cells.Item(11, 5).Value = 0.018540059432941663
$ws.Cells.Item(11, 6).Value = 0.08280530295777747
$ws.Cells.Item(11, 7).Value = -0.008121376914526635
$ws.Cells.Item(11, 8).Value = -0.0023361497597841523
$ws.Cells.Item(11, 9).Value = 0.02763002125397855
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = -0.025652894485768463

# Row 12: new row, date label + values
$ws.Cells.Item(12, 1).NumberFormat = "@"
$ws.Cells.Item(12, 1).Value = "2025-08-30"
$ws.Cells.Item(12, 1).Style = "Normal"
$ws.Cells.Item(12, 2).Value = 0.30955459402749685
$ws.Cells.Item(12, 3).Value = -0.1418878623126936
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(12, 5).Value = 0.0002577402581638719
$ws.Cells.Item(12, 6).Value = -0.00019365803327090083
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(12, 8).Value = -0.00025707184256496283
$ws.Cells.Item(12, 9).Value = -0.03523156852639539
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = -0.04484216455570389
